# Auto-generated script applying scheduled-runner price/profit updates
# to the Ultima_Profits leve-profitability sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1280.5769  # H40: was 1400.1428
$ws.Cells.Item(40, 9).Value = 1099  # I40: was 1137.5
$ws.Cells.Item(40, 10).Value = 1376.7059  # J40: was 1561.7693
$ws.Cells.Item(40, 11).Value = 1099  # K40: was 1137.5
$ws.Cells.Item(40, 12).Value = 1376.7059  # L40: was 1561.7693
$ws.Cells.Item(40, 13).Value = -924  # M40: was -962.5
$ws.Cells.Item(40, 14).Value = -1726.7059  # N40: was -1911.7693
$ws.Cells.Item(43, 8).Value = 6821.737  # H43: was 8084.4375
$ws.Cells.Item(43, 9).Value = 2800.375  # I43: was 2805.375
$ws.Cells.Item(43, 10).Value = 9746.362999999999  # J43: was 13363.5
$ws.Cells.Item(43, 11).Value = 2800.375  # K43: was 2805.375
$ws.Cells.Item(43, 12).Value = 9746.362999999999  # L43: was 13363.5
$ws.Cells.Item(43, 13).Value = -2731.375  # M43: was -2736.375
$ws.Cells.Item(43, 14).Value = -9884.362999999999  # N43: was -13501.5
$ws.Cells.Item(107, 8).Value = 5290.294  # H107: was 6075.6562
$ws.Cells.Item(107, 9).Value = 7023.75  # I107: was 8091.864
$ws.Cells.Item(107, 10).Value = 1130  # J107: was 1640
$ws.Cells.Item(107, 11).Value = 7023.75  # K107: was 8091.864
$ws.Cells.Item(107, 12).Value = 1130  # L107: was 1640
$ws.Cells.Item(107, 13).Value = -5103.75  # M107: was -6171.864
$ws.Cells.Item(107, 14).Value = -4970  # N107: was -5480
$ws.Cells.Item(132, 8).Value = 8261.6  # H132: was 6314.4443
$ws.Cells.Item(132, 9).Value = 5554.706  # I132: was 4427.9565
$ws.Cells.Item(132, 10).Value = 14013.75  # J132: was 9652.076999999999
$ws.Cells.Item(132, 11).Value = 16664.118  # K132: was 13283.8695
$ws.Cells.Item(132, 12).Value = 42041.25  # L132: was 28956.231
$ws.Cells.Item(132, 13).Value = -14134.118  # M132: was -10753.8695
$ws.Cells.Item(132, 14).Value = -47101.25  # N132: was -34016.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 15627831  # H61: was 16131906
$ws.Cells.Item(61, 9).Value = 20835920  # I61: was 21741764
$ws.Cells.Item(61, 11).Value = 20835920  # K61: was 21741764
$ws.Cells.Item(61, 13).Value = -20835708  # M61: was -21741552
$ws.Cells.Item(132, 8).Value = 8066296.5  # H132: was 8930482
$ws.Cells.Item(132, 9).Value = 9616590  # I132: was 10870854
$ws.Cells.Item(132, 11).Value = 28849770  # K132: was 32612562
$ws.Cells.Item(132, 13).Value = -28847240  # M132: was -32610032
$ws.Cells.Item(136, 8).Value = 15627831  # H136: was 16131906
$ws.Cells.Item(136, 9).Value = 20835920  # I136: was 21741764
$ws.Cells.Item(136, 11).Value = 62507760  # K136: was 65225292
$ws.Cells.Item(136, 13).Value = -62505210  # M136: was -65222742

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 25375.5  # H22: was 9497.909
$ws.Cells.Item(22, 9).Value = 500  # I22: was 447.8889
$ws.Cells.Item(22, 10).Value = 100002  # J22: was 50223
$ws.Cells.Item(22, 11).Value = 500  # K22: was 447.8889
$ws.Cells.Item(22, 12).Value = 100002  # L22: was 50223
$ws.Cells.Item(22, 13).Value = -327  # M22: was -274.8889
$ws.Cells.Item(22, 14).Value = -100348  # N22: was -50569
$ws.Cells.Item(36, 8).Value = 1561.6  # H36: was 3268.5
$ws.Cells.Item(36, 9).Value = 1561.6  # I36: was 3268.5
$ws.Cells.Item(36, 11).Value = 1561.6  # K36: was 3268.5
$ws.Cells.Item(36, 13).Value = -1027.6  # M36: was -2734.5
$ws.Cells.Item(105, 8).Value = 4693  # H105: was 4762
$ws.Cells.Item(105, 9).Value = 3666.1667  # I105: was 3799.6
$ws.Cells.Item(105, 10).Value = 4960.8696  # J105: was 4962.5
$ws.Cells.Item(105, 11).Value = 3666.1667  # K105: was 3799.6
$ws.Cells.Item(105, 12).Value = 4960.8696  # L105: was 4962.5
$ws.Cells.Item(105, 13).Value = -1919.1667  # M105: was -2052.6
$ws.Cells.Item(105, 14).Value = -8454.8696  # N105: was -8456.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1151.375  # H16: was 896.4783
$ws.Cells.Item(16, 9).Value = 862.4  # I16: was 771.0769
$ws.Cells.Item(16, 10).Value = 1633  # J16: was 1059.5
$ws.Cells.Item(16, 11).Value = 862.4  # K16: was 771.0769
$ws.Cells.Item(16, 12).Value = 1633  # L16: was 1059.5
$ws.Cells.Item(16, 13).Value = -575.4  # M16: was -484.0769
$ws.Cells.Item(16, 14).Value = -2207  # N16: was -1633.5
$ws.Cells.Item(99, 8).Value = 1239.909  # H99: was 1268.3
$ws.Cells.Item(99, 9).Value = 1094.6666  # I99: was 1122.4
$ws.Cells.Item(99, 11).Value = 1094.6666  # K99: was 1122.4
$ws.Cells.Item(99, 13).Value = 403.3334  # M99: was 375.5999999999999
$ws.Cells.Item(113, 8).Value = 1151.375  # H113: was 896.4783
$ws.Cells.Item(113, 9).Value = 862.4  # I113: was 771.0769
$ws.Cells.Item(113, 10).Value = 1633  # J113: was 1059.5
$ws.Cells.Item(113, 11).Value = 862.4  # K113: was 771.0769
$ws.Cells.Item(113, 12).Value = 1633  # L113: was 1059.5
$ws.Cells.Item(113, 13).Value = 1307.6  # M113: was 1398.9231
$ws.Cells.Item(113, 14).Value = -5973  # N113: was -5399.5
$ws.Cells.Item(122, 8).Value = 1537.6522  # H122: was 1454.1562
$ws.Cells.Item(122, 9).Value = 1461.3684  # I122: was 1409.0435
$ws.Cells.Item(122, 10).Value = 1900  # J122: was 1569.4445
$ws.Cells.Item(122, 11).Value = 4384.1052  # K122: was 4227.1305
$ws.Cells.Item(122, 12).Value = 5700  # L122: was 4708.333500000001
$ws.Cells.Item(122, 13).Value = -1934.1052  # M122: was -1777.1305
$ws.Cells.Item(122, 14).Value = -10600  # N122: was -9608.333500000001
$ws.Cells.Item(126, 8).Value = 1239.909  # H126: was 1268.3
$ws.Cells.Item(126, 9).Value = 1094.6666  # I126: was 1122.4
$ws.Cells.Item(126, 11).Value = 3283.9998  # K126: was 3367.2
$ws.Cells.Item(126, 13).Value = -813.9998000000001  # M126: was -897.2000000000003
$ws.Cells.Item(140, 8).Value = 35350  # H140: was 30442.857
$ws.Cells.Item(140, 10).Value = 35350  # J140: was 30442.857
$ws.Cells.Item(140, 12).Value = 35350  # L140: was 30442.857
$ws.Cells.Item(140, 14).Value = -45710  # N140: was -40802.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 1732.9166  # H121: was 1773
$ws.Cells.Item(121, 10).Value = 1732.9166  # J121: was 1773
$ws.Cells.Item(121, 12).Value = 5198.7498  # L121: was 5319
$ws.Cells.Item(121, 14).Value = -7818.7498  # N121: was -7939
$ws.Cells.Item(122, 8).Value = 861.55554  # H122: was 824.9259
$ws.Cells.Item(122, 9).Value = 767.95654  # I122: was 744.75
$ws.Cells.Item(122, 10).Value = 1399.75  # J122: was 1466.3334
$ws.Cells.Item(122, 11).Value = 6911.60886  # K122: was 6702.75
$ws.Cells.Item(122, 12).Value = 12597.75  # L122: was 13197.0006
$ws.Cells.Item(122, 13).Value = -4461.60886  # M122: was -4252.75
$ws.Cells.Item(122, 14).Value = -17497.75  # N122: was -18097.0006
$ws.Cells.Item(125, 10).Value = 5200  # J125: was 5062.5
$ws.Cells.Item(125, 12).Value = 15600  # L125: was 15187.5
$ws.Cells.Item(125, 14).Value = -25440  # N125: was -25027.5
$ws.Cells.Item(134, 8).Value = 4254.8  # H134: was 4311
$ws.Cells.Item(134, 9).Value = 2871.3333  # I134: was 3216.923
$ws.Cells.Item(134, 10).Value = 6330  # J134: was 6342.857
$ws.Cells.Item(134, 11).Value = 8613.999899999999  # K134: was 9650.769
$ws.Cells.Item(134, 12).Value = 18990  # L134: was 19028.571
$ws.Cells.Item(134, 13).Value = -3543.999899999999  # M134: was -4580.769
$ws.Cells.Item(134, 14).Value = -29130  # N134: was -29168.571
$ws.Cells.Item(137, 8).Value = 5754.1904  # H137: was 6177870.5
$ws.Cells.Item(137, 9).Value = 4392.375  # I137: was 18522424
$ws.Cells.Item(137, 10).Value = 6592.231  # J137: was 5594.3887
$ws.Cells.Item(137, 11).Value = 13177.125  # K137: was 55567272
$ws.Cells.Item(137, 12).Value = 19776.693  # L137: was 16783.1661
$ws.Cells.Item(137, 13).Value = -8077.125  # M137: was -55562172
$ws.Cells.Item(137, 14).Value = -29976.693  # N137: was -26983.1661

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 9714319  # H80: was 11535340
$ws.Cells.Item(80, 9).Value = 19610392  # I80: was 25643750
$ws.Cells.Item(80, 10).Value = 1703213  # J80: was 1882217
$ws.Cells.Item(80, 11).Value = 19610392  # K80: was 25643750
$ws.Cells.Item(80, 12).Value = 1703213  # L80: was 1882217
$ws.Cells.Item(80, 13).Value = -19609394  # M80: was -25642752
$ws.Cells.Item(80, 14).Value = -1705209  # N80: was -1884213
$ws.Cells.Item(83, 8).Value = 9714319  # H83: was 11535340
$ws.Cells.Item(83, 9).Value = 19610392  # I83: was 25643750
$ws.Cells.Item(83, 10).Value = 1703213  # J83: was 1882217
$ws.Cells.Item(83, 11).Value = 98051960  # K83: was 128218750
$ws.Cells.Item(83, 12).Value = 8516065  # L83: was 9411085
$ws.Cells.Item(83, 13).Value = -98046968  # M83: was -128213758
$ws.Cells.Item(83, 14).Value = -8526049  # N83: was -9421069

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5204.9546  # H7: was 4503.5415
$ws.Cells.Item(7, 9).Value = 5392.231  # I7: was 4363.2354
$ws.Cells.Item(7, 10).Value = 4934.4443  # J7: was 4844.2856
$ws.Cells.Item(7, 11).Value = 5392.231  # K7: was 4363.2354
$ws.Cells.Item(7, 12).Value = 4934.4443  # L7: was 4844.2856
$ws.Cells.Item(7, 13).Value = -5280.231  # M7: was -4251.2354
$ws.Cells.Item(7, 14).Value = -5158.4443  # N7: was -5068.2856
$ws.Cells.Item(22, 8).Value = 1440.3636  # H22: was 1080.1875
$ws.Cells.Item(22, 9).Value = 0  # I22: was 212.5
$ws.Cells.Item(22, 10).Value = 1440.3636  # J22: was 1369.4166
$ws.Cells.Item(22, 11).Value = 0  # K22: was 212.5
$ws.Cells.Item(22, 12).Value = 1440.3636  # L22: was 1369.4166
$ws.Cells.Item(22, 13).ClearContents()  # M22: was 82.5
$ws.Cells.Item(22, 14).Value = -2030.3636  # N22: was -1959.4166
$ws.Cells.Item(27, 8).Value = 1440.3636  # H27: was 1080.1875
$ws.Cells.Item(27, 9).Value = 0  # I27: was 212.5
$ws.Cells.Item(27, 10).Value = 1440.3636  # J27: was 1369.4166
$ws.Cells.Item(27, 11).Value = 0  # K27: was 212.5
$ws.Cells.Item(27, 12).Value = 1440.3636  # L27: was 1369.4166
$ws.Cells.Item(27, 13).ClearContents()  # M27: was -105.5
$ws.Cells.Item(27, 14).Value = -1654.3636  # N27: was -1583.4166
$ws.Cells.Item(40, 8).Value = 4637.625  # H40: was 4886.304
$ws.Cells.Item(40, 9).Value = 5484.385  # I40: was 5659.923
$ws.Cells.Item(40, 10).Value = 3636.9092  # J40: was 3880.6
$ws.Cells.Item(40, 11).Value = 5484.385  # K40: was 5659.923
$ws.Cells.Item(40, 12).Value = 3636.9092  # L40: was 3880.6
$ws.Cells.Item(40, 13).Value = -5348.385  # M40: was -5523.923
$ws.Cells.Item(40, 14).Value = -3908.9092  # N40: was -4152.6
$ws.Cells.Item(46, 8).Value = 1565.6666  # H46: was 1486.6666
$ws.Cells.Item(46, 9).Value = 1318.2  # I46: was 1375
$ws.Cells.Item(46, 10).Value = 1875  # J46: was 1614.2858
$ws.Cells.Item(46, 11).Value = 1318.2  # K46: was 1375
$ws.Cells.Item(46, 12).Value = 1875  # L46: was 1614.2858
$ws.Cells.Item(46, 13).Value = -1130.2  # M46: was -1187
$ws.Cells.Item(46, 14).Value = -2251  # N46: was -1990.2858
$ws.Cells.Item(61, 8).Value = 1275.4375  # H61: was 1159.5625
$ws.Cells.Item(61, 9).Value = 1009.4545  # I61: was 1039.5
$ws.Cells.Item(61, 10).Value = 1860.6  # J61: was 2000
$ws.Cells.Item(61, 11).Value = 1009.4545  # K61: was 1039.5
$ws.Cells.Item(61, 12).Value = 1860.6  # L61: was 2000
$ws.Cells.Item(61, 13).Value = -807.4545000000001  # M61: was -837.5
$ws.Cells.Item(61, 14).Value = -2264.6  # N61: was -2404
$ws.Cells.Item(101, 8).Value = 10362  # H101: was 15000
$ws.Cells.Item(101, 10).Value = 10362  # J101: was 15000
$ws.Cells.Item(101, 12).Value = 10362  # L101: was 15000
$ws.Cells.Item(101, 14).Value = -16852  # N101: was -21490
$ws.Cells.Item(113, 8).Value = 1275.4375  # H113: was 1159.5625
$ws.Cells.Item(113, 9).Value = 1009.4545  # I113: was 1039.5
$ws.Cells.Item(113, 10).Value = 1860.6  # J113: was 2000
$ws.Cells.Item(113, 11).Value = 1009.4545  # K113: was 1039.5
$ws.Cells.Item(113, 12).Value = 1860.6  # L113: was 2000
$ws.Cells.Item(113, 13).Value = 1160.5455  # M113: was 1130.5
$ws.Cells.Item(113, 14).Value = -6200.6  # N113: was -6340
$ws.Cells.Item(122, 8).Value = 6267.4165  # H122: was 5557.0312
$ws.Cells.Item(122, 9).Value = 6822.357  # I122: was 5930.294
$ws.Cells.Item(122, 10).Value = 5490.5  # J122: was 5134
$ws.Cells.Item(122, 11).Value = 20467.071  # K122: was 17790.882
$ws.Cells.Item(122, 12).Value = 16471.5  # L122: was 15402
$ws.Cells.Item(122, 13).Value = -18017.071  # M122: was -15340.882
$ws.Cells.Item(122, 14).Value = -21371.5  # N122: was -20302
$ws.Cells.Item(126, 8).Value = 5204.9546  # H126: was 4503.5415
$ws.Cells.Item(126, 9).Value = 5392.231  # I126: was 4363.2354
$ws.Cells.Item(126, 10).Value = 4934.4443  # J126: was 4844.2856
$ws.Cells.Item(126, 11).Value = 16176.693  # K126: was 13089.7062
$ws.Cells.Item(126, 12).Value = 14803.3329  # L126: was 14532.8568
$ws.Cells.Item(126, 13).Value = -13706.693  # M126: was -10619.7062
$ws.Cells.Item(126, 14).Value = -19743.3329  # N126: was -19472.8568
$ws.Cells.Item(132, 8).Value = 9440647  # H132: was 9265821
$ws.Cells.Item(132, 9).Value = 4464.68  # I132: was 4182.852
$ws.Cells.Item(132, 10).Value = 17865810  # J132: was 18527458
$ws.Cells.Item(132, 11).Value = 13394.04  # K132: was 12548.556
$ws.Cells.Item(132, 12).Value = 53597430  # L132: was 55582374
$ws.Cells.Item(132, 13).Value = -10864.04  # M132: was -10018.556
$ws.Cells.Item(132, 14).Value = -53602490  # N132: was -55587434

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 0  # H103: was 20301
$ws.Cells.Item(103, 10).Value = 0  # J103: was 20301
$ws.Cells.Item(103, 12).Value = 0  # L103: was 20301
$ws.Cells.Item(103, 14).ClearContents()  # N103: was -22645
